$d = $word.ActiveDocument

# 1. Add a new row to the "Requisitos funcionales" table that ends with RF_024
#    (table #4 of the document), containing RF_025 / Mostrar_imagenes_representativas.
#    Wrap the new text with a "_GoBack" bookmark, matching the target XML.
$t = $d.Tables.Item(4)
$newRow = $t.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "RF_025"

$cell2 = $newRow.Cells.Item(2)
$cell2.Range.Text = "Mostrar_imagenes_representativas"
$cellRange = $cell2.Range
$textOnly = $d.Range($cellRange.Start, $cellRange.End - 1)
$d.Bookmarks.Add("_GoBack", $textOnly)

# 2. Remove the old "_GoBack" bookmark that used to sit in the empty paragraph
#    right before the "Anexos" heading (it moved to the new table text above).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 3. Drop the stale lastRenderedPageBreak cached before the "Anexos" heading run;
#    re-writing the run text makes the engine recompute/drop the stale marker.
$d.Content.Find.Execute("Anexos", $true, $false, $false, $false, $false, $true, 1, $false, "Anexos", 2)

# 4. Update the footer's cached PAGE field result from 6 to 7 (new page count).
$sec = $d.Sections.Item(1)
$footer = $sec.Footers.Item(1)
$footer.Range.Find.Execute("6", $true, $false, $false, $false, $false, $true, 1, $false, "7", 2)
